# Audit-SEO.xlsx update:
#   - repurpose the old "black hat" row (row 16) into the
#     "Texte invisible page 2" entry
#   - add a new "accessiblité / couleur non conforme" row (17)
#   - add a new "accessiblité / police trop petite" row (18)
#   - extend column C formatting down through row 29
#   - move the active selection to C27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: rename the existing "black hat" entry --------------------
$ws.Range("B16").Value = "Texte invisible page 2"
$ws.Range("C16").Value = "texte non lisible par l’utilisateur"

# --- Row 17: new "couleur non conforme" entry --------------------------
$ws.Range("A17").Value = "accessiblité"
$ws.Range("A17").Font.Name = "Calibri"
$ws.Range("A17").Font.Size = 12
$ws.Range("A17").Font.Color = 0
$ws.Range("A17").Font.Family = 2
$ws.Range("A17").HorizontalAlignment = -4108

$ws.Range("B17").Value = "couleur non conforme"

$ws.Range("C17").Value = "couleur non conforme au norme  du niveau AA du WCAG 2.1"
$ws.Range("C17").Font.Name = "Arial"
$ws.Range("C17").Font.Size = 12
$ws.Range("C17").Font.Color = 0
$ws.Range("C17").HorizontalAlignment = -4108

$ws.Range("D17").Value = "correction des couleurs "

# --- Row 18: new "police trop petite" entry -----------------------------
$ws.Range("A18").Value = "accessiblité"
$ws.Range("A18").Font.Name = "Calibri"
$ws.Range("A18").Font.Size = 12
$ws.Range("A18").Font.Color = 0
$ws.Range("A18").Font.Family = 2
$ws.Range("A18").HorizontalAlignment = -4108

$ws.Range("B18").Value = "police trop petite "

$ws.Range("C18").Value = "police trop petite sur certain texte"
$ws.Range("C18").Font.Name = "Arial"
$ws.Range("C18").Font.Size = 12
$ws.Range("C18").Font.Color = 0
$ws.Range("C18").HorizontalAlignment = -4108

# --- Rows 19-29: extend column C formatting to match column B ----------
for ($r = 19; $r -le 29; $r++) {
    $ws.Cells.Item($r, 3).Font.Name = "Arial"
    $ws.Cells.Item($r, 3).Font.Size = 12
    $ws.Cells.Item($r, 3).Font.Color = 0
    $ws.Cells.Item($r, 3).HorizontalAlignment = -4108
}

# --- Selection moves to C27 --------------------------------------------
$ws.Range("C27").Select()
